$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.461.05"
$ws.Range("E2").Value = "  -4.84%  "
$ws.Range("D3").Value = "3.321.23"
$ws.Range("E3").Value = "  -5.70%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.00%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "3.309.61"
$ws.Range("E9").Value = "  -5.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.613"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.54%  "
$ws.Range("D15").Value = "3.852.01"
$ws.Range("E15").Value = "  -5.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.117"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.63%  "
$ws.Range("D17").Value = "3.315.08"
$ws.Range("E17").Value = "  -5.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.25%  "
$ws.Range("D20").Value = "63.331.47"
$ws.Range("E20").Value = "  -5.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.967"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "404.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.78%  "
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "575.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.17%  "
$ws.Range("E34").Value = "  -6.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.146"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "34.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("D40").Value = "0.0₃0732"
$ws.Range("E40").Value = "  -11.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.365"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.96%  "
$ws.Range("D42").Value = "3.144.90"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0401"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.127"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.49%  "
